{"js": "// The document has a handful of paragraphs whose <w:pPr> contains a stray,\n// overridden `Compact` paragraph style listed ahead of the paragraph's real\n// style (e.g. <w:pStyle w:val=\"Compact\"/><w:pStyle w:val=\"CasebookTitle\"/>).\n// Only the last <w:pStyle> is actually effective, so the \"Compact\" entry is\n// dead, invalid markup left over from an earlier export. Re-apply each\n// affected paragraph's own (already-effective) style so Word collapses the\n// paragraph properties down to the single, correct <w:pStyle> entry.\n\nconst targetStyles = [\n  \"Casebook Title\",\n  \"Casebook Subtitle\",\n  \"Section Number\",\n  \"Section Title\",\n  \"Section Subtitle\",\n  \"Resource Number\",\n  \"Resource Title\",\n  \"Resource Link\"\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/style\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  if (targetStyles.indexOf(paragraph.style) !== -1) {\n    // Re-assigning the same style forces Word to rewrite pPr with a single\n    // pStyle entry, dropping the stray leading \"Compact\" style.\n    paragraph.style = paragraph.style;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document has a handful of paragraphs whose <w:pPr> contains a stray,\n# overridden \"Compact\" paragraph style listed ahead of the paragraph's real\n# style (e.g. <w:pStyle w:val=\"Compact\"/><w:pStyle w:val=\"CasebookTitle\"/>).\n# Only the last <w:pStyle> is actually effective, so the \"Compact\" entry is\n# dead, invalid markup left over from an earlier export. Re-apply each\n# affected paragraph's own (already-effective) style so Word collapses the\n# paragraph properties down to the single, correct <w:pStyle> entry.\n\n$d = $word.ActiveDocument\n\n$targetStyles = @(\n    \"Casebook Title\",\n    \"Casebook Subtitle\",\n    \"Section Number\",\n    \"Section Title\",\n    \"Section Subtitle\",\n    \"Resource Number\",\n    \"Resource Title\",\n    \"Resource Link\"\n)\n\nforeach ($p in $d.Paragraphs) {\n    $styleName = $p.Style.NameLocal\n    if ($targetStyles -contains $styleName) {\n        # Re-assigning the same style forces Word to rewrite pPr with a\n        # single pStyle entry, dropping the stray leading \"Compact\" style.\n        $p.Style = $styleName\n    }\n}\n"}
